$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 1.196583684234042
$data[0,1] = 0.3420668026630835
$data[0,2] = 0
$data[0,3] = 0.4264942120488229
$data[0,4] = 0.4443680307746121
$data[0,5] = 0.002363637116791652
$data[0,6] = 0
$data[0,7] = 0.2114591347574857
$data[0,8] = 0
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 0.9672867250430812
$data[1,0] = 1.046102572727136
$data[1,1] = 0.3005729275224382
$data[1,2] = 0
$data[1,3] = 0.3719544669308021
$data[1,4] = 0.3878228170618172
$data[1,5] = 0.002366613124064439
$data[1,6] = 0
$data[1,7] = 0.219975127369783
$data[1,8] = 0
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 0.9845322939446532
$data[2,0] = 0.9533862914403812
$data[2,1] = 0.2749803947176304
$data[2,2] = 0
$data[2,3] = 0.3385623886600513
$data[2,4] = 0.3531389305169483
$data[2,5] = 0.002368533029871032
$data[2,6] = 0
$data[2,7] = 0.2255542313631516
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0.9967312379882287
$data[3,0] = 0.9155252024470997
$data[3,1] = 0.2645227655571887
$data[3,2] = 0
$data[3,3] = 0.3249769051013374
$data[3,4] = 0.3390132514313251
$data[3,5] = 0.002369338774571923
$data[3,6] = 0
$data[3,7] = 0.2279153481829637
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 1.002104328371558
$data[4,0] = 0.9092337212333064
$data[4,1] = 0.2627845795160795
$data[4,2] = 0
$data[4,3] = 0.3227223210947727
$data[4,4] = 0.336668177824194
$data[4,5] = 0.00236947398124443
$data[4,6] = 0
$data[4,7] = 0.2283126856130089
$data[4,8] = 0
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 1.003020714673468
$data[5,0] = 0.9528759984121393
$data[5,1] = 0.2748394740696085
$data[5,2] = 0
$data[5,3] = 0.3383790833323417
$data[5,4] = 0.3529483938344953
$data[5,5] = 0.002368543801739093
$data[5,6] = 0
$data[5,7] = 0.2255857202221101
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0.9968020779166835
$data[6,0] = 1.144765462767339
$data[6,1] = 0.3277838846451573
$data[6,2] = 0
$data[6,3] = 0.4076678744507518
$data[6,4] = 0.4248636149813478
$data[6,5] = 0.00236464406585063
$data[6,6] = 0
$data[6,7] = 0.214322462938366
$data[6,8] = 0
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0.9728970473210268
$data[7,0] = 1.518446807417149
$data[7,1] = 0.4306782330368151
$data[7,2] = 0
$data[7,3] = 0.5443964008934898
$data[7,4] = 0.5661985755041457
$data[7,5] = 0.002357728129746525
$data[7,6] = 0
$data[7,7] = 0.1950357710770465
$data[7,8] = 0
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 0.9389275659667931
$data[8,0] = 1.791329038949016
$data[8,1] = 0.5056943825953226
$data[8,2] = 0
$data[8,3] = 0.6455126250994425
$data[8,4] = 0.6702781546542269
$data[8,5] = 0.002353087950677465
$data[8,6] = 0
$data[8,7] = 0.1826020673085686
$data[8,8] = 0
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0.922022925739725
$data[9,0] = 1.915097322337829
$data[9,1] = 0.539692898105784
$data[9,2] = 0
$data[9,3] = 0.6916880816684454
$data[9,4] = 0.7176906081379002
$data[9,5] = 0.00235107170596921
$data[9,6] = 0
$data[9,7] = 0.1773288776496038
$data[9,8] = 0
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 0.916120723975439
$data[10,0] = 1.96191077613355
$data[10,1] = 0.5525486893321272
$data[10,2] = 0
$data[10,3] = 0.7092013862955184
$data[10,4] = 0.7356546913071611
$data[10,5] = 0.002350321729141085
$data[10,6] = 0
$data[10,7] = 0.1753876562856505
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 0.9141460393713601
$data[11,0] = 1.951831136435089
$data[11,1] = 0.549780804408158
$data[11,2] = 0
$data[11,3] = 0.7054283174838076
$data[11,4] = 0.7317853510981394
$data[11,5] = 0.00235048264931679
$data[11,6] = 0
$data[11,7] = 0.17580325047334
$data[11,8] = 0
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 0.914559691714544
$data[12,0] = 1.918949805726186
$data[12,1] = 0.5407509294969941
$data[12,2] = 0
$data[12,3] = 0.6931283447976426
$data[12,4] = 0.7191683204515869
$data[12,5] = 0.002351009734195058
$data[12,6] = 0
$data[12,7] = 0.1771680543044463
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 0.915953028912071
$data[13,0] = 1.89880182868734
$data[13,1] = 0.5352174196426063
$data[13,2] = 0
$data[13,3] = 0.6855979239808647
$data[13,4] = 0.7114413442032514
$data[13,5] = 0.002351334348869076
$data[13,6] = 0
$data[13,7] = 0.1780112959677265
$data[13,8] = 0
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 0.9168404907074716
$data[14,0] = 1.783232920515559
$data[14,1] = 0.5034699079803886
$data[14,2] = 0
$data[14,3] = 0.6424987017047243
$data[14,4] = 0.6671810134426437
$data[14,5] = 0.002353221614432716
$data[14,6] = 0
$data[14,7] = 0.1829544359565363
$data[14,8] = 0
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 0.922444902748154
$data[15,0] = 1.712239576930756
$data[15,1] = 0.4839610341875868
$data[15,2] = 0
$data[15,3] = 0.6161056492322814
$data[15,4] = 0.6400460337215605
$data[15,5] = 0.0023544035682763
$data[15,6] = 0
$data[15,7] = 0.1860853702093674
$data[15,8] = 0
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 0.9263432489847219
$data[16,0] = 1.671371642887777
$data[16,1] = 0.4727281481901855
$data[16,2] = 0
$data[16,3] = 0.6009416082069379
$data[16,4] = 0.6244449056556647
$data[16,5] = 0.00235509230490061
$data[16,6] = 0
$data[16,7] = 0.1879222141424739
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 0.9287535553534951
$data[17,0] = 1.657528612213014
$data[17,1] = 0.4689228576090727
$data[17,2] = 0
$data[17,3] = 0.5958100999789053
$data[17,4] = 0.6191636801734006
$data[17,5] = 0.002355327031343277
$data[17,6] = 0
$data[17,7] = 0.188550306498894
$data[17,8] = 0
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 0.9295984122913836
$data[18,0] = 1.719800519367197
$data[18,1] = 0.4860390220436557
$data[18,2] = 0
$data[18,3] = 0.6189135066490081
$data[18,4] = 0.642933953830422
$data[18,5] = 0.002354276825662348
$data[18,6] = 0
$data[18,7] = 0.1857483456903015
$data[18,8] = 0
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 0
$data[18,13] = 0.9259108440626136
$data[19,0] = 1.928609356399591
$data[19,1] = 0.5434037342341185
$data[19,2] = 0
$data[19,3] = 0.6967403766714995
$data[19,4] = 0.7228739723492197
$data[19,5] = 0.00235085454991631
$data[19,6] = 0
$data[19,7] = 0.1767656643904933
$data[19,8] = 0
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 0.9155366800606828
$data[20,0] = 2.064757001980013
$data[20,1] = 0.5807856220349095
$data[20,2] = 0
$data[20,3] = 0.7477670003336954
$data[20,4] = 0.7751780083420101
$data[20,5] = 0.002348696732730912
$data[20,6] = 0
$data[20,7] = 0.1712194226987487
$data[20,8] = 0
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 0.9102756130421312
$data[21,0] = 1.992122497827211
$data[21,1] = 0.5608443615416263
$data[21,2] = 0
$data[21,3] = 0.7205175549010079
$data[21,4] = 0.7472568307830727
$data[21,5] = 0.002349841210220663
$data[21,6] = 0
$data[21,7] = 0.17414968547107
$data[21,8] = 0
$data[21,9] = 0
$data[21,10] = 0
$data[21,11] = 0
$data[21,12] = 0
$data[21,13] = 0.9129434683832187
$data[22,0] = 1.716382381956635
$data[22,1] = 0.4850996165733932
$data[22,2] = 0
$data[22,3] = 0.6176440443557567
$data[22,4] = 0.6416283278902171
$data[22,5] = 0.002354334097266555
$data[22,6] = 0
$data[22,7] = 0.1859005998010996
$data[22,8] = 0
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 0
$data[22,13] = 0.9261058078616173
$data[23,0] = 1.417643088198815
$data[23,1] = 0.402943549143572
$data[23,2] = 0
$data[23,3] = 0.5073018126816606
$data[23,4] = 0.5279251897347166
$data[23,5] = 0.002359521280928687
$data[23,6] = 0
$data[23,7] = 0.199950531579419
$data[23,8] = 0
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 0
$data[23,13] = 0.9467163562435701

$ws.Range("B2:O25").Value = $data
Write-Host "Updated pl_mw values for 380 kV case (B2:O25)"
